# "update figure layout power plant"
#
# Figure 10 (power-plant layout drawing) relabels the transformer that sits
# between the medium-voltage collection grid and the high-voltage grid
# connection: "MV/LV transformer" -> "MV/HV transformer".
#
# That text lives in the TextBox named "TextBox 235" (shape id 236) on the
# 6th slide of the deck ("Power plant" figure), as the 4th top-level shape.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shape = $s.Shapes.Item(4)

if ($shape.Name -ne "TextBox 235") {
    # Defensive fallback: locate the shape by its current text instead of a
    # hard-coded index, in case shape ordering ever differs.
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $candidate = $s.Shapes.Item($i)
        if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText -and $candidate.TextFrame.TextRange.Text -eq "MV/LV transformer") {
            $shape = $candidate
            break
        }
    }
}

$shape.TextFrame.TextRange.Text = "MV/HV transformer"
